$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.489.26"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").Value = "'3.332.15"
$ws.Range("E3").Value = "  -0.48%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'587.28"
$ws.Range("E5").Value = "  +2.34%  "
$ws.Range("D6").Value = "'182.40"
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("D7").Value = "'0.651"
$ws.Range("E7").Value = "  +3.92%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'3.329.68"
$ws.Range("E9").Value = "  -0.56%  "
$ws.Range("E10").Value = "  -1.45%  "
$ws.Range("D11").Value = "'6.84"
$ws.Range("E11").Value = "  +2.25%  "
$ws.Range("E12").Value = "  +0.15%  "
$ws.Range("D13").Value = "'3.910.38"
$ws.Range("E13").Value = "  -0.52%  "
$ws.Range("E14").Value = "  -2.32%  "
$ws.Range("D15").Value = "'66.459.09"
$ws.Range("E15").Value = "  -0.39%  "
$ws.Range("D16").Value = "'26.57"
$ws.Range("E16").Value = "  -1.22%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "'3.371.52"
$ws.Range("E17").Value = "  +0.48%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "'0.0000166"
$ws.Range("E18").Value = "  -0.54%  "
$ws.Range("D19").Value = "'426.28"
$ws.Range("E19").Value = "  -2.23%  "
$ws.Range("E20").Value = "  -2.36%  "
$ws.Range("E21").Value = "  -3.20%  "
$ws.Range("D22").Value = "'7.40"
$ws.Range("E22").Value = "  -2.53%  "
$ws.Range("D23").Value = "'71.97"
$ws.Range("E23").Value = "  -2.09%  "
$ws.Range("E24").Value = "  +0.24%  "
$ws.Range("E25").Value = "  +0.27%  "
$ws.Range("D26").Value = "'3.465.83"
$ws.Range("E26").Value = "  -1.01%  "
$ws.Range("D27").Value = "'0.517"
$ws.Range("E27").Value = "  -0.36%  "
$ws.Range("D28").Value = "'0.204"
$ws.Range("E28").Value = "  +5.31%  "
$ws.Range("E29").Value = "  -1.22%  "
$ws.Range("D30").Value = "'9.04"
$ws.Range("E30").Value = "  -0.26%  "
$ws.Range("E31").Value = "  -0.13%  "
$ws.Range("D32").Value = "'1.93"
$ws.Range("E32").Value = "  -1.47%  "
$ws.Range("D33").Value = "'22.45"
$ws.Range("E33").Value = "  -1.79%  "
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("D35").Value = "'5.22"
$ws.Range("E35").Value = "  -1.24%  "
$ws.Range("D36").Value = "'6.65"
$ws.Range("E36").Value = "  -2.16%  "
$ws.Range("E37").Value = "  -2.54%  "
$ws.Range("D38").Value = "'160.89"
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("E39").Value = "  -2.23%  "
$ws.Range("E40").Value = "  +1.04%  "
$ws.Range("D41").Value = "'2.871.56"
$ws.Range("E41").Value = "  +1.97%  "
$ws.Range("D42").Value = "'26.51"
$ws.Range("E42").Value = "  -5.16%  "
$ws.Range("D43").Value = "'4.35"
$ws.Range("E43").Value = "  -2.13%  "
$ws.Range("D44").Value = "'0.761"
$ws.Range("E44").Value = "  -4.74%  "
$ws.Range("E45").Value = "  -1.40%  "
$ws.Range("E46").Value = "  -0.54%  "
$ws.Range("D47").Value = "'5.97"
$ws.Range("E47").Value = "  -3.78%  "
$ws.Range("E48").Value = "  -0.98%  "
$ws.Range("D49").Value = "'23.24"
$ws.Range("E49").Value = "  -4.40%  "
$ws.Range("D50").Value = "'314.22"
$ws.Range("E50").Value = "  -3.62%  "
$ws.Range("E51").Value = "  +0.19%  "
